# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# Concretely:
#   1. Insert a new "Player Info" worksheet (first tab) with the player's
#      ID / NAME / BATTING_HAND / BOWL_STYLE.
#   2. On "ODI Batting" and "ODI Bowling", rename the MATCH_CARD_LINK
#      column to MATCH_CODE and replace the full scorecard URL with just
#      the bare numeric match code that used to be the query-string tail.

$wb = $excel.ActiveWorkbook

# Helper: force a numeric-looking value to be written as *text*, matching
# the source workbook's convention of storing every cell (even numbers)
# as a string. Writing "'4634" makes Excel store it as text (quote
# prefix); re-applying the "Normal" style afterwards drops the transient
# quote-prefix formatting flag so the cell ends up as a plain text cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Insert a brand-new "Player Info" worksheet as the very first sheet.
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Font.Bold = $true

Set-TextValue $playerInfo.Range("A2") "5727"
$playerInfo.Range("B2").Value = "Naseem Shah"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# ---------------------------------------------------------------------
# 2. "ODI Batting": rename MATCH_CARD_LINK column to MATCH_CODE and
#    replace the URL values with the bare numeric match code.
# ---------------------------------------------------------------------
$oB = $wb.Worksheets.Item("ODI Batting")
$oB.Range("D1").Value = "MATCH_CODE"

$codes = @("4634", "4638", "4641", "4686", "4688")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $oB.Cells.Item($row, 4) $codes[$i]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": same rename/replace, column B this time.
# ---------------------------------------------------------------------
$oW = $wb.Worksheets.Item("ODI Bowling")
$oW.Range("B1").Value = "MATCH_CODE"

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    Set-TextValue $oW.Cells.Item($row, 2) $codes[$i]
}
